# Apply the commit: "Small adjustments to consumption. Added elysis SE03. Adjusted offshore wind SE."
#
# Workbook has two sheets: "Readme" (sheet 1) and "constants" (sheet 2).

$wb = $excel.ActiveWorkbook
$wsReadme = $wb.Worksheets.Item("Readme")
$wsConst  = $wb.Worksheets.Item("constants")

# ---------------------------------------------------------------------------
# 1. Readme sheet: remove the second paragraph of instructions (A3). The text
#    was removed from the shared string table entirely, so clearing the cell
#    contents (while keeping its style/formatting) reproduces that.
# ---------------------------------------------------------------------------
$wsReadme.Range("A3").ClearContents()

# ---------------------------------------------------------------------------
# 2. constants sheet: update numeric coefficients.
# ---------------------------------------------------------------------------

# Row 2: National Trends, 2025, C_0
$wsConst.Range("E2").Value = 80
$wsConst.Range("G2").Value = -600
$wsConst.Range("J2").Value = -600
$wsConst.Range("L2").Value = 0
$wsConst.Range("N2").Value = -705
$wsConst.Range("R2").Value = 140
$wsConst.Range("S2").Value = 130
$wsConst.Range("U2").Value = 75
$wsConst.Range("V2").Value = 65
$wsConst.Range("W2").Value = 600
$wsConst.Range("X2").Value = 0
$wsConst.Range("Y2").Value = 0
$wsConst.Range("Z2").Value = 0
$wsConst.Range("AA2").Value = 0

# Row 4: Distributed Energy, 2040, C_0
$wsConst.Range("E4").Value = 400
$wsConst.Range("J4").Value = 220
$wsConst.Range("K4").Value = 3000
$wsConst.Range("N4").Value = 1400
$wsConst.Range("R4").Value = 600
$wsConst.Range("T4").Value = 400
$wsConst.Range("V4").Value = 250
$wsConst.Range("Z4").Value = 1700

# Row 6: Distributed Energy, 2030, C_0
$wsConst.Range("E6").Value = 190
$wsConst.Range("G6").Value = 600
$wsConst.Range("J6").Value = 120
$wsConst.Range("K6").Value = 300
$wsConst.Range("N6").Value = 353
$wsConst.Range("R6").Value = 300
$wsConst.Range("T6").Value = 190
$wsConst.Range("W6").Value = 900
$wsConst.Range("Y6").Value = 150
$wsConst.Range("Z6").Value = 580
$wsConst.Range("AA6").Value = 150

# ---------------------------------------------------------------------------
# 3. Update view/selection state: "constants" becomes the active sheet,
#    selection moves to H11 there; "Readme" selection moves to A6 and is no
#    longer the selected tab.
# ---------------------------------------------------------------------------
$wsReadme.Range("A6").Select() | Out-Null
$wsConst.Select() | Out-Null
$wsConst.Range("H11").Select() | Out-Null

Write-Output "Edit applied."
